$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4 (slides/slide4.xml): new "graduation" title slide content
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$title4 = $s4.Shapes.Item(1).TextFrame.TextRange
$title4.Text = "毕业篇"

$sub4 = $s4.Shapes.Item(2).TextFrame.TextRange
# Setting straight to the final text can leave the original run-split
# ("--  " / old word) intact because of a shared prefix; break that by
# routing through an unrelated string first so the whole range collapses
# back down to a single run.
$sub4.Text = "~~placeholder~~"
$sub4.Text = "--  假期的两件事"

# ---------------------------------------------------------------------
# Slide 3 (slides/slide3.xml): author credit slide content
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Text = "荣明同学"
# Title textbox picks up "shrink text on overflow" autofit in the edited deck.
$s3.Shapes.Item(1).TextFrame.AutoSize = 2

$sub3 = $s3.Shapes.Item(2).TextFrame.TextRange
$sub3.Text = "~~placeholder~~"
$sub3.Text = "--  rmliu"

# Subtitle placeholder is repositioned to the left.
$subShape3 = $s3.Shapes.Item(2)
$subShape3.Left = 150.46251968503938
$subShape3.Top = 260.58748031496066

# ---------------------------------------------------------------------
# Touch the notes pane so the deck gains a notes master + (empty) notes
# slide, matching the authored edit.
# ---------------------------------------------------------------------
$notesPage = $p.Slides.Item(1).NotesPage
